# Insert a new data row at row 299 of the "Betarraga" sheet, shifting all
# subsequent rows (old 299..350) down by one (to 300..351), and populate the
# newly inserted row with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 299 (old row 299 -> becomes 300, etc.)
$ws.Rows(299).Insert()

# Populate the newly inserted row 299 with the new record's data.
$ws.Cells.Item(299, 1).Value = 11
$ws.Cells.Item(299, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(299, 3).Value = "Bíobío"
$ws.Cells.Item(299, 4).Value = 44776
$ws.Cells.Item(299, 5).Value = 8
$ws.Cells.Item(299, 6).Value = 100114014
$ws.Cells.Item(299, 7).Value = "Betarraga"
$ws.Cells.Item(299, 8).Value = "Sin especificar"
$ws.Cells.Item(299, 9).Value = "Primera"
$ws.Cells.Item(299, 10).Value = 300
$ws.Cells.Item(299, 11).Value = 600
$ws.Cells.Item(299, 12).Value = 650
$ws.Cells.Item(299, 13).Value = 625
$ws.Cells.Item(299, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(299, 15).Value = "Región Metropolitana"
$ws.Cells.Item(299, 16).Value = 125
$ws.Cells.Item(299, 17).Value = 5
$ws.Cells.Item(299, 18).Value = "Hortaliza"
